$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add new column E: "Mongoose Method" ---

# Header cell E2, matching the style of the other header cells (D2)
$ws.Range("E2").Value = "Mongoose Method"
$ws.Range("D2").Copy()
$ws.Range("E2").PasteSpecial(-4122)  # xlPasteFormats

# Data cells E3:E9 with the Mongoose method that corresponds to each route
$mongooseMethods = @(
    "Dog.find()",
    "N/A",
    "Dog.create()",
    "Dog.findById()",
    "Dog.findById",
    "Dog.findByIdAndUpdate()",
    "Dog.findByIdAndRemove()"
)

for ($i = 0; $i -lt $mongooseMethods.Length; $i++) {
    $row = 3 + $i
    $cell = $ws.Range("E$row")
    $cell.Value = $mongooseMethods[$i]
    $cell.Borders.LineStyle = 1
}

# Widen the new column to fit its contents
$ws.Columns.Item(5).ColumnWidth = 20.3

# Clear the previous clipboard marching-ants state
$excel.CutCopyMode = $false

# Update the active selection to reflect where editing left off
$ws.Range("D11").Select() | Out-Null
